# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Mon Jan  8 14:35:58 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.058.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "
$ws.Range("D3").Value = "'2.269.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.72%  "
$ws.Range("D5").Value = "'302.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "'95.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "'34.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.48%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("D12").Value = "'7.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "'2.610.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'2.267.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "'13.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "'0.802"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.59%  "
$ws.Range("D18").Value = "'44.929.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").Value = "'12.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.69%  "
$ws.Range("D20").Value = "'0.0₃0924"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").Value = "'65.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'238.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -4.21%  "
$ws.Range("D27").Value = "'41.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.51%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("D29").Value = "'9.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").Value = "'19.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "'153.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.36%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("E34").Value = "  -2.28%  "
$ws.Range("D35").Value = "'2.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "'0.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.16%  "
$ws.Range("D38").Value = "'1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0312"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("D41").Value = "'3.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").Value = "'13.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.50%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").Value = "'1.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.92%  "
$ws.Range("D45").Value = "'1.750.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.50%  "
$ws.Range("E46").Value = "  +3.40%  "
$ws.Range("D47").Value = "'71.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").Value = "'76.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.78%  "
$ws.Range("D49").Value = "'96.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'54.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").Value = "'4.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.42%  "
